# Apply MV -datos- update: revise Q1 2021 figures and append Q2 2021 (01-04-2021) row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 75 ("01-01-2021") with revised figures ---
$ws.Range("B75").Value = 1577946
$ws.Range("D75").Value = 199347
$ws.Range("E75").Value = 118035
$ws.Range("F75").Value = 81312
$ws.Range("G75").Value = 222158
$ws.Range("I75").Value = 193808
$ws.Range("J75").Value = 686
$ws.Range("K75").Value = 244318
$ws.Range("M75").Value = 206074
$ws.Range("N75").Value = 546380
$ws.Range("O75").Value = 406543
$ws.Range("V75").Value = 161076
$ws.Range("W75").Value = -19579
$ws.Range("X75").Value = 1597525
$ws.Range("AB75").Value = 235642
$ws.Range("AD75").Value = 210369
$ws.Range("AE75").Value = 306084
$ws.Range("AF75").Value = 41721
$ws.Range("AG75").Value = 264364
$ws.Range("AH75").Value = 513642
$ws.Range("AI75").Value = 436350
$ws.Range("AP75").Value = 157124

# --- Append new row 76 for period "01-04-2021" ---
$a76 = $ws.Cells.Item(76, 1)
# Force the date-looking label to be stored as text, matching the other "Serie" labels,
# then restore the default (unstyled) cell appearance used by the rest of column A.
$a76.NumberFormat = "@"
$a76.Value = "01-04-2021"
$a76.Style = "Normal"

$ws.Range("B76").Value = 1588167
$ws.Range("C76").Value = 439
$ws.Range("D76").Value = 214604
$ws.Range("E76").Value = 130414
$ws.Range("F76").Value = 84190
$ws.Range("G76").Value = 222189
$ws.Range("H76").Value = 35712
$ws.Range("I76").Value = 186477
$ws.Range("J76").Value = 552
$ws.Range("K76").Value = 254803
$ws.Range("L76").Value = 45361
$ws.Range("M76").Value = 209441
$ws.Range("N76").Value = 538975
$ws.Range("O76").Value = 400968
$ws.Range("P76").Value = 26896
$ws.Range("Q76").Value = 111111
$ws.Range("R76").Value = 147978
$ws.Range("S76").Value = 43693
$ws.Range("T76").Value = 42184
$ws.Range("U76").Value = 1509
$ws.Range("V76").Value = 164933
$ws.Range("W76").Value = -12280
$ws.Range("X76").Value = 1600447
$ws.Range("Y76").Value = 192716
$ws.Range("Z76").Value = 119399
$ws.Range("AA76").Value = 73318
$ws.Range("AB76").Value = 239044
$ws.Range("AC76").Value = 35307
$ws.Range("AD76").Value = 203737
$ws.Range("AE76").Value = 315252
$ws.Range("AF76").Value = 49160
$ws.Range("AG76").Value = 266092
$ws.Range("AH76").Value = 501538
$ws.Range("AI76").Value = 425630
$ws.Range("AJ76").Value = 26038
$ws.Range("AK76").Value = 49871
$ws.Range("AL76").Value = 147978
$ws.Range("AM76").Value = 43693
$ws.Range("AN76").Value = 42184
$ws.Range("AO76").Value = 1509
$ws.Range("AP76").Value = 160224
